# planilha_gincana_solidaria.xlsx — "Add files via upload" re-edit
#
# 1) doacoes_registros!H2:H15 had a hard-coded default "5" in the
#    Pontos_Unit column; clear those back to blank input cells and give
#    the whole H2:H68 column the same yellow "fill in me" look as the
#    neighbouring D/F/G input columns (copy format from G2).
# 2) Selection/active-sheet bookkeeping: participantes moves its cursor to
#    G22 and is no longer the active tab; doacoes_registros becomes the
#    active tab with the cursor at G74.

$wb = $excel.ActiveWorkbook

$participantes = $wb.Worksheets.Item("participantes")
$doacoes = $wb.Worksheets.Item("doacoes_registros")

# --- doacoes_registros: blank out the pre-filled "5" defaults and match
#     the input-cell formatting used by the rest of the row ---
[void]$doacoes.Range("G2").Copy()
[void]$doacoes.Range("H2:H68").PasteSpecial(-4122)  # xlPasteFormats
[void]$doacoes.Range("H2:H15").ClearContents()

# --- selection / active sheet bookkeeping ---
# Leave participantes' cursor on G22 without it being the active tab.
[void]$participantes.Range("G22").Select()

# doacoes_registros becomes the active tab, cursor parked at G74.
[void]$doacoes.Activate()
[void]$doacoes.Range("G74").Select()
